{"js": "const replacements = [\n  [\"2025-05-04 Sunday\", \"2025-05-05 Monday\"],\n  [\"22\u00d763=\", \"33\u00d765=\"],\n  [\"18\u00d765=\", \"86\u00d725=\"],\n  [\"12\u00d720=\", \"92\u00d787=\"],\n  [\"95\u00d738=\", \"87\u00d793=\"],\n  [\"29\u00d713=\", \"97\u00d733=\"],\n  [\"72\u00d715=\", \"17\u00d794=\"],\n  [\"25\u00d754=\", \"55\u00d794=\"],\n  [\"12\u00d778=\", \"59\u00d717=\"],\n  [\"21\u00d726=\", \"57\u00d749=\"],\n  [\"84\u00d771=\", \"70\u00d742=\"],\n  [\"54\u00d732=\", \"73\u00d753=\"],\n  [\"57\u00d790=\", \"92\u00d762=\"],\n  [\"36\u00d717=\", \"32\u00d713=\"],\n  [\"81\u00d723=\", \"11\u00d770=\"],\n  [\"61\u00d770=\", \"73\u00d765=\"],\n  [\"26\u00d786=\", \"12\u00d795=\"],\n  [\"92\u00d794=\", \"79\u00d782=\"],\n  [\"33\u00d756=\", \"31\u00d783=\"],\n  [\"40\u00d728=\", \"47\u00d795=\"],\n  [\"65\u00d798=\", \"55\u00d744=\"],\n  [\"44\u00d769=\", \"82\u00d761=\"],\n  [\"18\u00d754=\", \"67\u00d794=\"],\n  [\"12\u00d728=\", \"19\u00d784=\"],\n  [\"89\u00d712=\", \"50\u00d799=\"],\n  [\"57\u00d762=\", \"16\u00d771=\"],\n];\n\nfor (const [oldText, newText] of replacements) {\n  const results = context.document.body.search(oldText, { matchCase: true, matchWholeWord: false });\n  results.load(\"items\");\n  await context.sync();\n  for (const r of results.items) {\n    r.insertText(newText, \"Replace\");\n  }\n  await context.sync();\n}\n", "ps1": "$d = $word.ActiveDocument\n\n$replacements = @(\n    @(\"2025-05-04 Sunday\", \"2025-05-05 Monday\"),\n    @(\"22\u00d763=\", \"33\u00d765=\"),\n    @(\"18\u00d765=\", \"86\u00d725=\"),\n    @(\"12\u00d720=\", \"92\u00d787=\"),\n    @(\"95\u00d738=\", \"87\u00d793=\"),\n    @(\"29\u00d713=\", \"97\u00d733=\"),\n    @(\"72\u00d715=\", \"17\u00d794=\"),\n    @(\"25\u00d754=\", \"55\u00d794=\"),\n    @(\"12\u00d778=\", \"59\u00d717=\"),\n    @(\"21\u00d726=\", \"57\u00d749=\"),\n    @(\"84\u00d771=\", \"70\u00d742=\"),\n    @(\"54\u00d732=\", \"73\u00d753=\"),\n    @(\"57\u00d790=\", \"92\u00d762=\"),\n    @(\"36\u00d717=\", \"32\u00d713=\"),\n    @(\"81\u00d723=\", \"11\u00d770=\"),\n    @(\"61\u00d770=\", \"73\u00d765=\"),\n    @(\"26\u00d786=\", \"12\u00d795=\"),\n    @(\"92\u00d794=\", \"79\u00d782=\"),\n    @(\"33\u00d756=\", \"31\u00d783=\"),\n    @(\"40\u00d728=\", \"47\u00d795=\"),\n    @(\"65\u00d798=\", \"55\u00d744=\"),\n    @(\"44\u00d769=\", \"82\u00d761=\"),\n    @(\"18\u00d754=\", \"67\u00d794=\"),\n    @(\"12\u00d728=\", \"19\u00d784=\"),\n    @(\"89\u00d712=\", \"50\u00d799=\"),\n    @(\"57\u00d762=\", \"16\u00d771=\")\n)\n\nforeach ($pair in $replacements) {\n    $find = $pair[0]\n    $replace = $pair[1]\n    $range = $d.Content\n    $range.Find.ClearFormatting()\n    $range.Find.Replacement.ClearFormatting()\n    $range.Find.Text = $find\n    $range.Find.Replacement.Text = $replace\n    $range.Find.Forward = $true\n    $range.Find.Wrap = 1\n    $range.Find.MatchCase = $true\n    $range.Find.MatchWholeWord = $false\n    $range.Find.MatchWildcards = $false\n    $range.Find.Execute($find, $true, $false, $false, $false, $false, $true, 1, $false, $replace, 2)\n}\n"}
